# Applies the "GUI Kameras / Cursor Kameras / Update ManagementSummary" doc edits.
#
# Helper: replace the first occurrence of $oldSub inside a Paragraph/TextRange
# object with $newSub. Using Characters(start,len) (rather than re-assigning
# the whole paragraph .Text) makes the host split only the run(s) touching the
# edited span, leaving every other run (and its rPr/formatting) untouched -
# mirroring what PowerPoint itself does on an in-place edit.
function Replace-InRange($range, $oldSub, $newSub) {
    $full = $range.Text
    # Paragraph/TextRange .Text carries a trailing CR (paragraph mark) in
    # this host; strip it before searching so indices/lengths stay in bounds.
    if ($full.Length -gt 0 -and [int][char]$full[$full.Length - 1] -eq 13) {
        $full = $full.Substring(0, $full.Length - 1)
    }
    $idx = $full.IndexOf($oldSub)
    if ($idx -lt 0) {
        throw "Replace-InRange: substring not found: [$oldSub]"
    }
    $sub = $range.Characters($idx + 1, $oldSub.Length)
    $sub.Text = $newSub
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape "Textfeld 8" (Ergebnis box): GUI wording + run-merge tidy-up.
# ---------------------------------------------------------------------
$ergebnisShape = $s.Shapes.Item(4)
$ergebnisText = $ergebnisShape.TextFrame.TextRange
$coveringPara = $ergebnisText.Paragraphs(3)

# "...Darstellung der UI-Elemente auf eine..." -> "...GUI-Elemente..."
Replace-InRange $coveringPara "UI-Elemente " "GUI-Elemente "

# Merge the "... Wand, dem " / "primären Inputgerät..." runs back into one
# contiguous run (same text, same formatting). Rewriting the whole tail of
# the paragraph (from the start of " Wand, dem " through the end) as a
# single assignment collapses it into one run instead of the original two.
$fullPara = $coveringPara.Text
if ($fullPara.Length -gt 0 -and [int][char]$fullPara[$fullPara.Length - 1] -eq 13) {
    $fullPara = $fullPara.Substring(0, $fullPara.Length - 1)
}
$mergeIdx = $fullPara.IndexOf(" Wand, dem ")
$mergeTail = $fullPara.Substring($mergeIdx)
$mergeRange = $coveringPara.Characters($mergeIdx + 1, $mergeTail.Length)
$mergeRange.Text = $mergeTail

# ---------------------------------------------------------------------
# Shape "Textfeld 12" (Ausgangslage und Umsetzung box).
# ---------------------------------------------------------------------
$umsetzungShape = $s.Shapes.Item(7)

# Grow the text box to match the expanded copy below.
$umsetzungShape.Height = 419.25472440944884

$umsetzungText = $umsetzungShape.TextFrame.TextRange
$hauptPara = $umsetzungText.Paragraphs(5)

Replace-InRange $hauptPara "war die " "ist die "
Replace-InRange $hauptPara "wurde die " "wird die "
Replace-InRange $hauptPara "für die Manipulationen " "für sämtliche Manipulationen "
